# Apply "nuevos experimentos no convexos" data update.
# Updates recorded experiment values (gamma/restriction coefficients, the
# modified point, vec_bf and vec_BF) on several sheets. Every one of these
# cells stores its value as text (shared string) in the workbook, even
# though most of them look like plain numbers, so plain numeric-looking
# strings must be written with a leading apostrophe (forcing Excel to treat
# them as literal text) and then have their cell style restored to "Normal"
# so no left-over text formatting is applied to the cell.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )

    if ($Text -match '^-?\d+(\.\d+)?$') {
        # Numeric-looking text: prefix with an apostrophe so Excel keeps it
        # as text instead of converting it to a real number, then restore
        # the default "Normal" style so the quote-prefix formatting that
        # Excel applies to the cell is not left behind.
        $Range.Value = "'" + $Text
        $Range.Style = "Normal"
    }
    else {
        $Range.Value = $Text
    }
}

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsFollower.Range("A2") "8.95 - y"
Set-TextValue $wsFollower.Range("B2") "-8.95"

Set-TextValue $wsFollower.Range("D2") "0.68"
Set-TextValue $wsFollower.Range("E2") "1.4000000000000001"
Set-TextValue $wsFollower.Range("F2") "7.800000000000001"

Set-TextValue $wsFollower.Range("A3") "-1.9499999999999993 - x + y"
Set-TextValue $wsFollower.Range("B3") "-1.0500000000000007"

Set-TextValue $wsFollower.Range("D3") "0.24"
Set-TextValue $wsFollower.Range("E3") "10.0"
Set-TextValue $wsFollower.Range("F3") "5.5"

Set-TextValue $wsFollower.Range("A4") "-24.9 + x + 2y"
Set-TextValue $wsFollower.Range("B4") "12.899999999999999"

Set-TextValue $wsFollower.Range("D4") "0.44"
Set-TextValue $wsFollower.Range("E4") "0.5"
Set-TextValue $wsFollower.Range("F4") "0.8999999999999999"

Set-TextValue $wsFollower.Range("A5") "-19.98 + 4x - y"
Set-TextValue $wsFollower.Range("B5") "7.050000000000001"

Set-TextValue $wsFollower.Range("D5") "1.0"
Set-TextValue $wsFollower.Range("E5") "2.9"
Set-TextValue $wsFollower.Range("F5") "2.4"

$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "7.0"
Set-TextValue $wsPunto.Range("B2") "8.95"

$wsVecBf = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $wsVecBf.Range("A2") "-0.43999999999999995"

$wsVecBF = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $wsVecBF.Range("A2") "-1.0999999999999996"
Set-TextValue $wsVecBF.Range("A3") "-3.6999999999999993"
